$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, scraped from the "cryptos" refresh diff.
# Values that look numeric (e.g. "28.79", "0.110", "1.00") must be forced to
# Text so Excel does not silently coerce them to numbers and drop meaningful
# trailing/leading zeros or the dotted-thousands "65.434.81" style.
$updates = [ordered]@{
    "D2" = '65.434.81'
    "E2" = '  +3.73%  '
    "D3" = '3.492.97'
    "E3" = '  +2.97%  '
    "E4" = '  +0.00%  '
    "D5" = '581.03'
    "D6" = '162.92'
    "E6" = '  +5.01%  '
    "D7" = '0.612'
    "E7" = '  +12.39%  '
    "E8" = '  +0.00%  '
    "D9" = '3.496.50'
    "E9" = '  +3.03%  '
    "E10" = '  -1.68%  '
    "D11" = '0.126'
    "E11" = '  +3.78%  '
    "D12" = '0.448'
    "E12" = '  +3.75%  '
    "D13" = '4.099.38'
    "E13" = '  +3.08%  '
    "D14" = '0.135'
    "E14" = '  +0.57%  '
    "D15" = '0.0000194'
    "E15" = '  +2.93%  '
    "D16" = '28.79'
    "D17" = '65.416.38'
    "E17" = '  +3.58%  '
    "D18" = '3.486.75'
    "E18" = '  +1.41%  '
    "D19" = '6.47'
    "E19" = '  +3.56%  '
    "D20" = '14.41'
    "E20" = '  +2.44%  '
    "D21" = '386.40'
    "E21" = '  +2.60%  '
    "D22" = '8.25'
    "E22" = '  +2.61%  '
    "D23" = '0.553'
    "E23" = '  +4.60%  '
    "D24" = '72.87'
    "E24" = '  +2.19%  '
    "E25" = '  +0.26%  '
    "D26" = '0.0000121'
    "E26" = '  +3.24%  '
    "D27" = '10.15'
    "E27" = '  +7.97%  '
    "E28" = '  +0.90%  '
    "D29" = '1.00'
    "E29" = '  -0.11%  '
    "E30" = '  +13.01%  '
    "D31" = '6.22'
    "E31" = '  +2.29%  '
    "E32" = '  +3.51%  '
    "D33" = '23.78'
    "E33" = '  +2.80%  '
    "D34" = '7.20'
    "E34" = '  +6.34%  '
    "D35" = '1.63'
    "E35" = '  +12.23%  '
    "D36" = '162.92'
    "E36" = '  +2.00%  '
    "D37" = '1.94'
    "E37" = '  +6.41%  '
    "D38" = '3.034.62'
    "E38" = '  +2.36%  '
    "D39" = '0.0781'
    "E39" = '  +3.16%  '
    "D40" = '27.08'
    "E40" = '  +0.25%  '
    "B41" = 'RenderToken'
    "C41" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "D41" = '6.83'
    "E41" = '  +6.44%  '
    "B42" = 'Filecoin'
    "C42" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D42" = '4.62'
    "E42" = '  +7.17%  '
    "D43" = '0.0322'
    "E43" = '  +1.18%  '
    "D44" = '43.06'
    "E44" = '  +3.47%  '
    "D45" = '0.784'
    "E45" = '  +3.84%  '
    "D46" = '25.96'
    "E46" = '  +11.45%  '
    "D47" = '1.12'
    "E47" = '  +5.04%  '
    "D48" = '320.78'
    "E48" = '  +10.96%  '
    "B49" = 'Cosmos'
    "C49" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "D49" = '6.77'
    "E49" = '  +6.56%  '
    "B50" = 'SuiNetwork'
    "C50" = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    "D50" = '0.884'
    "E50" = '  +6.36%  '
    "B51" = 'Stellar'
    "C51" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "D51" = '0.110'
    "E51" = '  +6.85%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($value.Trim() -match "^[+-]?[0-9]*\.?[0-9]+$") {
        # Pre-format as Text so the numeric-looking string is kept verbatim,
        # then drop the format override again so the cell style matches the
        # rest of the (unstyled) data column once the value is committed.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}
